# WRI Updates to HK model from 11/15
# - CESTR sheet: rename "Tax Rate (dimensionless)" header to "Tax Rate"
#   and zero out the tax rate value (per July 2019 workshop: set 0 for Hong Kong)
# - About sheet: add a red note explaining the change, and make it the
#   newly selected cell; make CESTR the active/selected sheet.

$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")
$cestr = $wb.Worksheets.Item("CESTR")

# --- CESTR sheet edits ---
$cestr.Range("B1").Value = "Tax Rate"
$cestr.Range("B2").Value = 0

# --- About sheet edits: new explanatory note in red font ---
$about.Range("B13").Value = "according to July 2019 workshop, set 0 to hong kong"
$about.Range("B13").Font.Color = 255

# --- Selection / active sheet bookkeeping to match the saved view state ---
[void]$about.Range("B13").Select()
[void]$cestr.Activate()
[void]$cestr.Range("B2").Select()
